# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Choclo, Choclero - Primera/Segunda,
# Region de O'Higgins, Fecha 44567) at the top of the existing data block,
# pushing the rest of the historical rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("46:47").Insert()

# New row 46: Choclero / Primera
$ws.Cells.Item(46,1).Value  = 11
$ws.Cells.Item(46,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(46,3).Value  = "Bíobío"
$ws.Cells.Item(46,4).Value  = 44567
$ws.Cells.Item(46,5).Value  = 8
$ws.Cells.Item(46,6).Value  = 100112024
$ws.Cells.Item(46,7).Value  = "Choclo"
$ws.Cells.Item(46,8).Value  = "Choclero"
$ws.Cells.Item(46,9).Value  = "Primera"
$ws.Cells.Item(46,10).Value = 20000
$ws.Cells.Item(46,11).Value = 250
$ws.Cells.Item(46,12).Value = 300
$ws.Cells.Item(46,13).Value = 275
$ws.Cells.Item(46,14).Value = "`$/unidad"
$ws.Cells.Item(46,15).Value = "Región de O'Higgins"
$ws.Cells.Item(46,16).Value = 275
$ws.Cells.Item(46,17).Value = 1
$ws.Cells.Item(46,18).Value = "Hortaliza"

# New row 47: Choclero / Segunda
$ws.Cells.Item(47,1).Value  = 11
$ws.Cells.Item(47,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(47,3).Value  = "Bíobío"
$ws.Cells.Item(47,4).Value  = 44567
$ws.Cells.Item(47,5).Value  = 8
$ws.Cells.Item(47,6).Value  = 100112024
$ws.Cells.Item(47,7).Value  = "Choclo"
$ws.Cells.Item(47,8).Value  = "Choclero"
$ws.Cells.Item(47,9).Value  = "Segunda"
$ws.Cells.Item(47,10).Value = 10000
$ws.Cells.Item(47,11).Value = 200
$ws.Cells.Item(47,12).Value = 200
$ws.Cells.Item(47,13).Value = 200
$ws.Cells.Item(47,14).Value = "`$/unidad"
$ws.Cells.Item(47,15).Value = "Región de O'Higgins"
$ws.Cells.Item(47,16).Value = 200
$ws.Cells.Item(47,17).Value = 1
$ws.Cells.Item(47,18).Value = "Hortaliza"
